$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update data values (row 5) ---
$ws.Range("C5").Value = 1620
$ws.Range("D5").Value = 1600
$ws.Range("G5").Value = 2070

# --- Update data values (row 6) ---
$ws.Range("G6").Value = 2100

# --- Update data values (row 8) ---
$ws.Range("D8").Value = 1320

# --- Update data values (row 9) ---
$ws.Range("D9").Value = 1250
$ws.Range("G9").Value = 800

$wb.Application.Calculate()

# --- Update sheet view (scroll position / selection) ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C13").Select()
